# Update Il1a-Il1r2 LR-pairs sheet with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"7.134879666666667"
$ws.Cells.Item(2, 8).Value = [double]"21.404639"
$ws.Cells.Item(2, 9).Value = [double]"0.07716103050836744"
$ws.Cells.Item(2, 10).Value = [double]"0.07716103050836744"
$ws.Cells.Item(2, 13).Value = [double]"0.301183"
$ws.Cells.Item(2, 14).Value = [double]"0.602366"
$ws.Cells.Item(2, 15).Value = [double]"0.0004546669609434378"
$ws.Cells.Item(2, 16).Value = [double]"0.0003032119723243698"
$ws.Cells.Item(2, 17).Value = [double]"2.148904462645667"
$ws.Cells.Item(2, 18).Value = [double]"12.893426775874"
$ws.Cells.Item(2, 19).Value = [double]"3.508257124450331E-05"
$ws.Cells.Item(2, 20).Value = [double]"2.339614824702296E-05"
$ws.Cells.Item(3, 7).Value = [double]"7.134879666666667"
$ws.Cells.Item(3, 8).Value = [double]"21.404639"
$ws.Cells.Item(3, 9).Value = [double]"0.07716103050836744"
$ws.Cells.Item(3, 10).Value = [double]"0.07716103050836744"
$ws.Cells.Item(3, 15).Value = [double]"0.002672378908982388"
$ws.Cells.Item(3, 16).Value = [double]"0.002673266421566756"
$ws.Cells.Item(3, 17).Value = [double]"12.630534999676"
$ws.Cells.Item(3, 18).Value = [double]"113.674814997084"
$ws.Cells.Item(3, 19).Value = [double]"0.0002062035105259078"
$ws.Cells.Item(3, 20).Value = [double]"0.0002062719919115067"
$ws.Cells.Item(4, 7).Value = [double]"7.134879666666667"
$ws.Cells.Item(4, 8).Value = [double]"21.404639"
$ws.Cells.Item(4, 9).Value = [double]"0.07716103050836744"
$ws.Cells.Item(4, 10).Value = [double]"0.07716103050836744"
$ws.Cells.Item(4, 13).Value = [double]"24.13741566666667"
$ws.Cells.Item(4, 14).Value = [double]"72.41224700000001"
$ws.Cells.Item(4, 15).Value = [double]"0.03643793117869155"
$ws.Cells.Item(4, 16).Value = [double]"0.03645003242764271"
$ws.Cells.Item(4, 17).Value = [double]"172.2175562459815"
$ws.Cells.Item(4, 18).Value = [double]"1549.958006213833"
$ws.Cells.Item(4, 19).Value = [double]"0.002811588319340812"
$ws.Cells.Item(4, 20).Value = [double]"0.002812522064180322"
$ws.Cells.Item(5, 7).Value = [double]"7.134879666666667"
$ws.Cells.Item(5, 8).Value = [double]"21.404639"
$ws.Cells.Item(5, 9).Value = [double]"0.07716103050836744"
$ws.Cells.Item(5, 10).Value = [double]"0.07716103050836744"
$ws.Cells.Item(5, 13).Value = [double]"0.358584"
$ws.Cells.Item(5, 14).Value = [double]"0.717168"
$ws.Cells.Item(5, 15).Value = [double]"0.0005413197209767541"
$ws.Cells.Item(5, 16).Value = [double]"0.000360999664270433"
$ws.Cells.Item(5, 17).Value = [double]"2.558453690392"
$ws.Cells.Item(5, 18).Value = [double]"15.350722142352"
$ws.Cells.Item(5, 19).Value = [double]"4.176878750506827E-05"
$ws.Cells.Item(5, 20).Value = [double]"2.785510610828128E-05"
$ws.Cells.Item(6, 7).Value = [double]"7.134879666666667"
$ws.Cells.Item(6, 8).Value = [double]"21.404639"
$ws.Cells.Item(6, 9).Value = [double]"0.07716103050836744"
$ws.Cells.Item(6, 10).Value = [double]"0.07716103050836744"
$ws.Cells.Item(6, 13).Value = [double]"629.608429"
$ws.Cells.Item(6, 14).Value = [double]"1888.825287"
$ws.Cells.Item(6, 15).Value = [double]"0.9504591925766137"
$ws.Cells.Item(6, 16).Value = [double]"0.9507748456045226"
$ws.Cells.Item(6, 17).Value = [double]"4492.180378034043"
$ws.Cells.Item(6, 18).Value = [double]"40429.6234023064"
$ws.Cells.Item(6, 19).Value = [double]"0.07333841075536238"
$ws.Cells.Item(6, 20).Value = [double]"0.07336276686827892"
$ws.Cells.Item(7, 7).Value = [double]"7.134879666666667"
$ws.Cells.Item(7, 8).Value = [double]"21.404639"
$ws.Cells.Item(7, 9).Value = [double]"0.07716103050836744"
$ws.Cells.Item(7, 10).Value = [double]"0.07716103050836744"
$ws.Cells.Item(7, 13).Value = [double]"6.249660666666666"
$ws.Cells.Item(7, 14).Value = [double]"18.748982"
$ws.Cells.Item(7, 15).Value = [double]"0.009434510653792123"
$ws.Cells.Item(7, 16).Value = [double]"0.009437643909673034"
$ws.Cells.Item(7, 17).Value = [double]"44.59057681416644"
$ws.Cells.Item(7, 18).Value = [double]"401.3151913274979"
$ws.Cells.Item(7, 19).Value = [double]"0.0007279765643887716"
$ws.Cells.Item(7, 20).Value = [double]"0.0007282183296413891"
$ws.Cells.Item(8, 5).Value = [double]"3"
$ws.Cells.Item(8, 6).Value = [double]"1"
$ws.Cells.Item(8, 7).Value = [double]"81.06813666666666"
$ws.Cells.Item(8, 8).Value = [double]"243.20441"
$ws.Cells.Item(8, 9).Value = [double]"0.8767212985829614"
$ws.Cells.Item(8, 10).Value = [double]"0.8767212985829616"
$ws.Cells.Item(8, 13).Value = [double]"0.301183"
$ws.Cells.Item(8, 14).Value = [double]"0.602366"
$ws.Cells.Item(8, 15).Value = [double]"0.0004546669609434378"
$ws.Cells.Item(8, 16).Value = [double]"0.0003032119723243698"
$ws.Cells.Item(8, 17).Value = [double]"24.41634460567666"
$ws.Cells.Item(8, 18).Value = [double]"146.49806763406"
$ws.Cells.Item(8, 19).Value = [double]"0.0003986162084210994"
$ws.Cells.Item(8, 20).Value = [double]"0.0002658323941221225"
$ws.Cells.Item(9, 5).Value = [double]"3"
$ws.Cells.Item(9, 6).Value = [double]"1"
$ws.Cells.Item(9, 7).Value = [double]"81.06813666666666"
$ws.Cells.Item(9, 8).Value = [double]"243.20441"
$ws.Cells.Item(9, 9).Value = [double]"0.8767212985829614"
$ws.Cells.Item(9, 10).Value = [double]"0.8767212985829616"
$ws.Cells.Item(9, 15).Value = [double]"0.002672378908982388"
$ws.Cells.Item(9, 16).Value = [double]"0.002673266421566756"
$ws.Cells.Item(9, 17).Value = [double]"143.51103107044"
$ws.Cells.Item(9, 18).Value = [double]"1291.59927963396"
$ws.Cells.Item(9, 19).Value = [double]"0.002342931507388757"
$ws.Cells.Item(9, 20).Value = [double]"0.002343709608574233"
$ws.Cells.Item(10, 5).Value = [double]"3"
$ws.Cells.Item(10, 6).Value = [double]"1"
$ws.Cells.Item(10, 7).Value = [double]"81.06813666666666"
$ws.Cells.Item(10, 8).Value = [double]"243.20441"
$ws.Cells.Item(10, 9).Value = [double]"0.8767212985829614"
$ws.Cells.Item(10, 10).Value = [double]"0.8767212985829616"
$ws.Cells.Item(10, 13).Value = [double]"24.13741566666667"
$ws.Cells.Item(10, 14).Value = [double]"72.41224700000001"
$ws.Cells.Item(10, 15).Value = [double]"0.03643793117869155"
$ws.Cells.Item(10, 16).Value = [double]"0.03645003242764271"
$ws.Cells.Item(10, 17).Value = [double]"1956.775312045474"
$ws.Cells.Item(10, 18).Value = [double]"17610.97780840927"
$ws.Cells.Item(10, 19).Value = [double]"0.03194591034065904"
$ws.Cells.Item(10, 20).Value = [double]"0.03195651976335398"
$ws.Cells.Item(11, 5).Value = [double]"3"
$ws.Cells.Item(11, 6).Value = [double]"1"
$ws.Cells.Item(11, 7).Value = [double]"81.06813666666666"
$ws.Cells.Item(11, 8).Value = [double]"243.20441"
$ws.Cells.Item(11, 9).Value = [double]"0.8767212985829614"
$ws.Cells.Item(11, 10).Value = [double]"0.8767212985829616"
$ws.Cells.Item(11, 13).Value = [double]"0.358584"
$ws.Cells.Item(11, 14).Value = [double]"0.717168"
$ws.Cells.Item(11, 15).Value = [double]"0.0005413197209767541"
$ws.Cells.Item(11, 16).Value = [double]"0.000360999664270433"
$ws.Cells.Item(11, 17).Value = [double]"29.06973671848"
$ws.Cells.Item(11, 18).Value = [double]"174.41842031088"
$ws.Cells.Item(11, 19).Value = [double]"0.0004745865287233062"
$ws.Cells.Item(11, 20).Value = [double]"0.0003164960944471872"
$ws.Cells.Item(12, 5).Value = [double]"3"
$ws.Cells.Item(12, 6).Value = [double]"1"
$ws.Cells.Item(12, 7).Value = [double]"81.06813666666666"
$ws.Cells.Item(12, 8).Value = [double]"243.20441"
$ws.Cells.Item(12, 9).Value = [double]"0.8767212985829614"
$ws.Cells.Item(12, 10).Value = [double]"0.8767212985829616"
$ws.Cells.Item(12, 13).Value = [double]"629.608429"
$ws.Cells.Item(12, 14).Value = [double]"1888.825287"
$ws.Cells.Item(12, 15).Value = [double]"0.9504591925766137"
$ws.Cells.Item(12, 16).Value = [double]"0.9507748456045226"
$ws.Cells.Item(12, 17).Value = [double]"51041.18216865729"
$ws.Cells.Item(12, 18).Value = [double]"459370.6395179157"
$ws.Cells.Item(12, 19).Value = [double]"0.8332878175658818"
$ws.Cells.Item(12, 20).Value = [double]"0.8335645572984118"
$ws.Cells.Item(13, 5).Value = [double]"3"
$ws.Cells.Item(13, 6).Value = [double]"1"
$ws.Cells.Item(13, 7).Value = [double]"81.06813666666666"
$ws.Cells.Item(13, 8).Value = [double]"243.20441"
$ws.Cells.Item(13, 9).Value = [double]"0.8767212985829614"
$ws.Cells.Item(13, 10).Value = [double]"0.8767212985829616"
$ws.Cells.Item(13, 13).Value = [double]"6.249660666666666"
$ws.Cells.Item(13, 14).Value = [double]"18.748982"
$ws.Cells.Item(13, 15).Value = [double]"0.009434510653792123"
$ws.Cells.Item(13, 16).Value = [double]"0.009437643909673034"
$ws.Cells.Item(13, 17).Value = [double]"506.6483450456244"
$ws.Cells.Item(13, 18).Value = [double]"4559.83510541062"
$ws.Cells.Item(13, 19).Value = [double]"0.008271436431887415"
$ws.Cells.Item(13, 20).Value = [double]"0.008274183424052121"
$ws.Cells.Item(14, 7).Value = [double]"4.264381"
$ws.Cells.Item(14, 8).Value = [double]"12.793143"
$ws.Cells.Item(14, 9).Value = [double]"0.04611767090867112"
$ws.Cells.Item(14, 10).Value = [double]"0.04611767090867113"
$ws.Cells.Item(14, 13).Value = [double]"0.301183"
$ws.Cells.Item(14, 14).Value = [double]"0.602366"
$ws.Cells.Item(14, 15).Value = [double]"0.0004546669609434378"
$ws.Cells.Item(14, 16).Value = [double]"0.0003032119723243698"
$ws.Cells.Item(14, 17).Value = [double]"1.284359062723"
$ws.Cells.Item(14, 18).Value = [double]"7.706154376338"
$ws.Cells.Item(14, 19).Value = [double]"2.096818127783509E-05"
$ws.Cells.Item(14, 20).Value = [double]"1.398342995522438E-05"
$ws.Cells.Item(15, 7).Value = [double]"4.264381"
$ws.Cells.Item(15, 8).Value = [double]"12.793143"
$ws.Cells.Item(15, 9).Value = [double]"0.04611767090867112"
$ws.Cells.Item(15, 10).Value = [double]"0.04611767090867113"
$ws.Cells.Item(15, 15).Value = [double]"0.002672378908982388"
$ws.Cells.Item(15, 16).Value = [double]"0.002673266421566756"
$ws.Cells.Item(15, 17).Value = [double]"7.549028994012"
$ws.Cells.Item(15, 18).Value = [double]"67.941260946108"
$ws.Cells.Item(15, 19).Value = [double]"0.0001232438910677234"
$ws.Cells.Item(15, 20).Value = [double]"0.0001232848210810166"
$ws.Cells.Item(16, 7).Value = [double]"4.264381"
$ws.Cells.Item(16, 8).Value = [double]"12.793143"
$ws.Cells.Item(16, 9).Value = [double]"0.04611767090867112"
$ws.Cells.Item(16, 10).Value = [double]"0.04611767090867113"
$ws.Cells.Item(16, 13).Value = [double]"24.13741566666667"
$ws.Cells.Item(16, 14).Value = [double]"72.41224700000001"
$ws.Cells.Item(16, 15).Value = [double]"0.03643793117869155"
$ws.Cells.Item(16, 16).Value = [double]"0.03645003242764271"
$ws.Cells.Item(16, 17).Value = [double]"102.9311367580357"
$ws.Cells.Item(16, 18).Value = [double]"926.3802308223211"
$ws.Cells.Item(16, 19).Value = [double]"0.001680432518691704"
$ws.Cells.Item(16, 20).Value = [double]"0.001680990600108418"
$ws.Cells.Item(17, 7).Value = [double]"4.264381"
$ws.Cells.Item(17, 8).Value = [double]"12.793143"
$ws.Cells.Item(17, 9).Value = [double]"0.04611767090867112"
$ws.Cells.Item(17, 10).Value = [double]"0.04611767090867113"
$ws.Cells.Item(17, 13).Value = [double]"0.358584"
$ws.Cells.Item(17, 14).Value = [double]"0.717168"
$ws.Cells.Item(17, 15).Value = [double]"0.0005413197209767541"
$ws.Cells.Item(17, 16).Value = [double]"0.000360999664270433"
$ws.Cells.Item(17, 17).Value = [double]"1.529138796504"
$ws.Cells.Item(17, 18).Value = [double]"9.174832779024001"
$ws.Cells.Item(17, 19).Value = [double]"2.496440474837962E-05"
$ws.Cells.Item(17, 20).Value = [double]"1.664846371496459E-05"
$ws.Cells.Item(18, 7).Value = [double]"4.264381"
$ws.Cells.Item(18, 8).Value = [double]"12.793143"
$ws.Cells.Item(18, 9).Value = [double]"0.04611767090867112"
$ws.Cells.Item(18, 10).Value = [double]"0.04611767090867113"
$ws.Cells.Item(18, 13).Value = [double]"629.608429"
$ws.Cells.Item(18, 14).Value = [double]"1888.825287"
$ws.Cells.Item(18, 15).Value = [double]"0.9504591925766137"
$ws.Cells.Item(18, 16).Value = [double]"0.9507748456045226"
$ws.Cells.Item(18, 17).Value = [double]"2684.890222067449"
$ws.Cells.Item(18, 18).Value = [double]"24164.01199860704"
$ws.Cells.Item(18, 19).Value = [double]"0.04383296425536955"
$ws.Cells.Item(18, 20).Value = [double]"0.04384752143783198"
$ws.Cells.Item(19, 7).Value = [double]"4.264381"
$ws.Cells.Item(19, 8).Value = [double]"12.793143"
$ws.Cells.Item(19, 9).Value = [double]"0.04611767090867112"
$ws.Cells.Item(19, 10).Value = [double]"0.04611767090867113"
$ws.Cells.Item(19, 13).Value = [double]"6.249660666666666"
$ws.Cells.Item(19, 14).Value = [double]"18.748982"
$ws.Cells.Item(19, 15).Value = [double]"0.009434510653792123"
$ws.Cells.Item(19, 16).Value = [double]"0.009437643909673034"
$ws.Cells.Item(19, 17).Value = [double]"26.65093420338066"
$ws.Cells.Item(19, 18).Value = [double]"239.858407830426"
$ws.Cells.Item(19, 19).Value = [double]"0.0004350976575159368"
$ws.Cells.Item(19, 20).Value = [double]"0.0004352421559795253"
